# Fix sale quantity prediction: amend result into integral value.
# The predicted sale-quantity values stored in column C for rows 2241-2515
# (the tail of the forecast series) are truncated down to whole numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2241
$lastRow  = 2515
$col      = 3   # column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value()
    $cell.Value = [int]$current
}

# Update the view to reflect where the edit was made: scrolled so row 2219
# is at the top, with the amended range selected.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 2219
$win.ScrollColumn = 1

$ws.Range("C" + $firstRow + ":C" + $lastRow).Select()
